$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the report title (row 1) - date changed from 2025-09-16 to 2025-09-17
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = '萊爾富 工作統計表  篩選月份：202509   (  製表日期:2025-09-17  )'

# ---------------------------------------------------------------------------
# 2. Fix formatting on existing row 89 (P89 / AC89 now wrap text)
# ---------------------------------------------------------------------------
$ws.Range("P89").WrapText = $true
$ws.Range("AC89").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Add new row 90 (copy format from row 88, then set values)
# ---------------------------------------------------------------------------
$ws.Range("A88:AK88").Copy()
$ws.Range("A90:AK90").PasteSpecial(-4122)
$ws.Range("P90").WrapText = $true
$ws.Range("AC90").WrapText = $true

$ws.Cells.Item(90, 1).Value = 88
$ws.Cells.Item(90, 2).Value = '維修'
$ws.Cells.Item(90, 3).Value = 2025092221
$ws.Cells.Item(90, 4).Value = 'E4282114091601'
$ws.Cells.Item(90, 5).Value = '一般件'
$ws.Cells.Item(90, 6).Value = 4282
$ws.Cells.Item(90, 7).Value = '林口建林店'
$ws.Cells.Item(90, 8).Value = '新北市林口區'
$ws.Cells.Item(90, 9).Value = '2025-09-16 16:20:36'
$ws.Cells.Item(90, 10).Value = '星期二'
$ws.Cells.Item(90, 11).Value = '下午'
$ws.Cells.Item(90, 12).Value = 'HL59'
$ws.Cells.Item(90, 13).Value = 'HL-LIFE-ET熱感機T70II'
$ws.Cells.Item(90, 14).Value = 5902
$ws.Cells.Item(90, 15).Value = '印字不清'
$ws.Cells.Item(90, 16).Value = '門市反映MMK熱感機(T70-II)列印出的單據到兩台TM都無法刷讀，門市已有重裝紙捲重開機仍異常...須請台芝到店協助(沒辦法操作)'
$ws.Cells.Item(90, 17).Value = 'THILF04282'
$ws.Cells.Item(90, 18).Value = '新北一'
$ws.Cells.Item(90, 19).Value = '湯家瑋'
$ws.Cells.Item(90, 20).Value = 1
$ws.Cells.Item(90, 21).Value = '已完工'
$ws.Cells.Item(90, 22).Value = '2025-09-16 16:34:11'
$ws.Cells.Item(90, 23).Value = '2025-09-17 13:31:00'
$ws.Cells.Item(90, 24).Value = '2025-09-17 14:30:00'
$ws.Cells.Item(90, 25).Value = '2025-09-17 20:34:00'
$ws.Cells.Item(90, 26).Value = 1
$ws.Cells.Item(90, 28).Value = '到場處理'
$ws.Cells.Item(90, 29).Value = '清潔t70'
$ws.Cells.Item(90, 37).Value = 'O'

# ---------------------------------------------------------------------------
# 4. Add new row 91 (copy format from the now-updated row 89, then set values)
# ---------------------------------------------------------------------------
$ws.Range("A89:AK89").Copy()
$ws.Range("A91:AK91").PasteSpecial(-4122)

$ws.Cells.Item(91, 1).Value = 89
$ws.Cells.Item(91, 2).Value = '服務'
$ws.Cells.Item(91, 3).Value = 2025092284
$ws.Cells.Item(91, 6).Value = 3999
$ws.Cells.Item(91, 7).Value = '新莊福祐店'
$ws.Cells.Item(91, 8).Value = '新北市新莊區'
$ws.Cells.Item(91, 17).Value = 'THILF03999'
$ws.Cells.Item(91, 18).Value = '新北一'
$ws.Cells.Item(91, 19).Value = '湯家瑋'
$ws.Cells.Item(91, 20).Value = 1
$ws.Cells.Item(91, 21).Value = '已完工'
$ws.Cells.Item(91, 22).Value = '2025-09-17 10:40:36'
$ws.Cells.Item(91, 23).Value = '2025-09-17 10:10:00'
$ws.Cells.Item(91, 24).Value = '2025-09-17 10:40:00'
$ws.Cells.Item(91, 26).Value = 0.5
$ws.Cells.Item(91, 28).Value = '到場處理'
$ws.Cells.Item(91, 29).Value = 'PMQ3'
$ws.Cells.Item(91, 30).Value = 'O'
$ws.Cells.Item(91, 37).Value = 'O'

# ---------------------------------------------------------------------------
# 5. Add new row 92 (copy format from row 88, fix P92/AC92 to non-wrap, set values)
# ---------------------------------------------------------------------------
$ws.Range("A88:AK88").Copy()
$ws.Range("A92:AK92").PasteSpecial(-4122)
$ws.Range("M88").Copy()
$ws.Range("P92").PasteSpecial(-4122)
$ws.Range("AC92").PasteSpecial(-4122)

$ws.Cells.Item(92, 1).Value = 90
$ws.Cells.Item(92, 2).Value = '服務'
$ws.Cells.Item(92, 3).Value = 2025092314
$ws.Cells.Item(92, 6).Value = 3929
$ws.Cells.Item(92, 7).Value = '蘆洲中山一'
$ws.Cells.Item(92, 8).Value = '新北市蘆洲區'
$ws.Cells.Item(92, 17).Value = 'THILF03929'
$ws.Cells.Item(92, 18).Value = '新北一'
$ws.Cells.Item(92, 19).Value = '吳宗鴻'
$ws.Cells.Item(92, 20).Value = 1
$ws.Cells.Item(92, 21).Value = '已完工'
$ws.Cells.Item(92, 22).Value = '2025-09-17 13:42:12'
$ws.Cells.Item(92, 23).Value = '2025-09-17 13:25:00'
$ws.Cells.Item(92, 24).Value = '2025-09-17 13:50:00'
$ws.Cells.Item(92, 26).Value = 0.4
$ws.Cells.Item(92, 29).Value = '裝潢回裝完成'
$ws.Cells.Item(92, 31).Value = 'O'
$ws.Cells.Item(92, 37).Value = 'O'

# ---------------------------------------------------------------------------
# 6. Update the print area to cover the new rows
# ---------------------------------------------------------------------------
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "Report!Print_Area") {
        $dn.RefersTo = "='Report'!`$A`$1:`$AK`$92"
    }
}

# ---------------------------------------------------------------------------
# 7. Restore default gridlines / headers display and set the new selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true
$ws.Range("AC89").Select()
